# InterpolatedResults.xlsx — fix position errors in csvParser.m output.
#
# 1) Header row (row 1): rename the shared-string headers from the raw
#    Mechanical-export identifiers to human readable labels, and rename
#    the first column from "NodeNumbers" to "CutLocation".
# 2) Column A (rows 2-103): replace the stale "1" placeholder with the
#    actual cut-location values (3.00, 3.02, 3.04, ... in 0.02" steps).
# 3) Column D (rows 2-103): these were mistakenly populated with the
#    Z-location numeric series; reset them to the literal value 1.
# 4) A handful of column widths shift slightly to fit the new header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header text -----------------------------------------------------
$ws.Range("A1").Value = "CutLocation"
$ws.Range("B1").Value = "X Locations (inches)"
$ws.Range("C1").Value = "Y Locations (inches)"
$ws.Range("D1").Value = "Z Locations (inches)"
$ws.Range("E1").Value = "Equivalent Elastic Strains (in/in)"
$ws.Range("F1").Value = "Equivalent Stress (psi)"
$ws.Range("G1").Value = "Max Principal Elastic Strain (in/in)"
$ws.Range("H1").Value = "Max Principal Stress (psi)"
$ws.Range("I1").Value = "Middle Principal Elastic Strain (in/in)"
$ws.Range("J1").Value = "Middle Principal Stress (psi)"
$ws.Range("K1").Value = "Min Principal Elastic Strain (in/in)"
$ws.Range("L1").Value = "Min Principal Stress (psi)"
$ws.Range("M1").Value = "Total Deformation (in)"

# --- 2) Column A: correct cut-location values (rows 2..103) ------------
$cutLocations = @(3, 3.02, 3.04, 3.0600000000000001, 3.0800000000000001, 3.1000000000000001, 3.1200000000000001, 3.1400000000000001, 3.1600000000000001, 3.1800000000000002, 3.2000000000000002, 3.2200000000000002, 3.2400000000000002, 3.2599999999999998, 3.2799999999999998, 3.2999999999999998, 3.3199999999999998, 3.3399999999999999, 3.3599999999999999, 3.3799999999999999, 3.3999999999999999, 3.4199999999999999, 3.4399999999999999, 3.46, 3.48, 3.5, 3.52, 3.54, 3.5600000000000001, 3.5800000000000001, 3.6000000000000001, 3.6200000000000001, 3.6400000000000001, 3.6600000000000001, 3.6800000000000002, 3.7000000000000002, 3.7200000000000002, 3.7400000000000002, 3.7599999999999998, 3.7799999999999998, 3.7999999999999998, 3.8199999999999998, 3.8399999999999999, 3.8599999999999999, 3.8799999999999999, 3.8999999999999999, 3.9199999999999999, 3.9399999999999999, 3.96, 3.98, 4, 4.0199999999999996, 4.04, 4.0599999999999996, 4.0800000000000001, 4.0999999999999996, 4.1200000000000001, 4.1399999999999997, 4.1600000000000001, 4.1799999999999997, 4.2000000000000002, 4.2199999999999998, 4.2400000000000002, 4.2599999999999998, 4.2800000000000002, 4.2999999999999998, 4.3200000000000003, 4.3399999999999999, 4.3600000000000003, 4.3799999999999999, 4.4000000000000004, 4.4199999999999999, 4.4400000000000004, 4.46, 4.4800000000000004, 4.5, 4.5199999999999996, 4.54, 4.5599999999999996, 4.5800000000000001, 4.5999999999999996, 4.6200000000000001, 4.6399999999999997, 4.6600000000000001, 4.6799999999999997, 4.7000000000000002, 4.7199999999999998, 4.7400000000000002, 4.7599999999999998, 4.7800000000000002, 4.7999999999999998, 4.8200000000000003, 4.8399999999999999, 4.8600000000000003, 4.8799999999999999, 4.9000000000000004, 4.9199999999999999, 4.9400000000000004, 4.96, 4.9800000000000004, 5, 5.0199999999999996)

for ($i = 0; $i -lt $cutLocations.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $cutLocations[$i]
}

# --- 3) Column D: reset mis-populated values to 1 (rows 2..103) --------
for ($row = 2; $row -le 103; $row++) {
    $ws.Cells.Item($row, 4).Value = 1
}

# --- 4) Column width tweaks to fit the new header text ------------------
# (target OOXML widths: A=11.046875, E=26.37890625, G=28.046875,
#  H=21.046875, I=30.046875, J=23.046875, K=27.6015625, L=20.6015625;
#  ColumnWidth values below are the closest this host's 1/6-character
#  width quantization can reach.)
$ws.Columns.Item(1).ColumnWidth = 10.166666666666666
$ws.Columns.Item(5).ColumnWidth = 25.5
$ws.Columns.Item(7).ColumnWidth = 27.166666666666668
$ws.Columns.Item(8).ColumnWidth = 20.166666666666668
$ws.Columns.Item(9).ColumnWidth = 29.166666666666668
$ws.Columns.Item(10).ColumnWidth = 22.166666666666668
$ws.Columns.Item(11).ColumnWidth = 26.833333333333332
$ws.Columns.Item(12).ColumnWidth = 19.833333333333332

Write-Host "done"
